$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 23:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1585860
$ws.Range("C4").Value = 15277
$ws.Range("D4").Value = 366591
$ws.Range("E4").Value = 1124649
$ws.Range("G4").Value = 1087
$ws.Range("H4").Value = 94620

# Row 11 - Alemania
$ws.Range("B11").Value = 178486
$ws.Range("C11").Value = 659
$ws.Range("E11").Value = 13321

# Row 60 - Kazajistan
$ws.Range("D60").Value = 3734
$ws.Range("E60").Value = 3200

# Row 83 - Costa de Marfil
$ws.Range("B83").Value = 2231
$ws.Range("C83").Value = 78
$ws.Range("D83").Value = 1083
$ws.Range("E83").Value = 1119
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 29

# Row 153 - Birmania
$ws.Range("B153").Value = 199
$ws.Range("C153").Value = 6
$ws.Range("D153").Value = 108
$ws.Range("E153").Value = 85

# Row 180 - Zimbabue
$ws.Range("B180").Value = 48
$ws.Range("C180").Value = 2
$ws.Range("E180").Value = 26
